$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PANDUAN")
$rng = $ws.Range("M3:M1048576")
$rng.ClearFormats()
